$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "wahid"
$ws.Range("A3").Value = "abass"

$ws.Range("B2").Value = 9442571997
$ws.Range("C2").Value = 35

$ws.Range("B3").Value = 5632562356

$ws.Range("C2").Select()
